$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.162.75"
$ws.Range("E2").Value = "  -6.78%  "

$ws.Range("D3").Value = "'2.891.90"
$ws.Range("E3").Value = "  -9.98%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'528.79"
$ws.Range("E5").Value = "  -10.99%  "

$ws.Range("D6").Value = "'126.62"
$ws.Range("E6").Value = "  -16.51%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").Value = "'2.870.85"
$ws.Range("E8").Value = "  -10.40%  "

$ws.Range("D9").Value = "'0.449"
$ws.Range("E9").Value = "  -17.83%  "

$ws.Range("D10").Value = "'0.139"
$ws.Range("E10").Value = "  -20.00%  "

$ws.Range("E11").Value = "  -14.24%  "

$ws.Range("D12").Value = "'0.416"
$ws.Range("E12").Value = "  -16.38%  "

$ws.Range("D13").Value = "'30.73"
$ws.Range("E13").Value = "  -21.54%  "

$ws.Range("D14").Value = "'0.0000193"
$ws.Range("E14").Value = "  -20.79%  "

$ws.Range("D15").Value = "'3.358.33"
$ws.Range("E15").Value = "  -10.08%  "

$ws.Range("D16").Value = "'61.985.28"
$ws.Range("E16").Value = "  -7.13%  "

$ws.Range("E17").Value = "  -5.83%  "

$ws.Range("D18").Value = "'2.884.99"
$ws.Range("E18").Value = "  -10.37%  "

$ws.Range("D19").Value = "'460.97"
$ws.Range("E19").Value = "  -13.65%  "

$ws.Range("D20").Value = "'5.97"
$ws.Range("E20").Value = "  -16.69%  "

$ws.Range("D21").Value = "'12.29"
$ws.Range("E21").Value = "  -18.17%  "

$ws.Range("D22").Value = "'0.610"
$ws.Range("E22").Value = "  -19.82%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'11.66"
$ws.Range("E23").Value = "  -15.84%  "

$ws.Range("D24").Value = "'72.89"
$ws.Range("E24").Value = "  -15.34%  "

$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("B26").Value = "Uniswap"
$ws.Range("C26").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D26").Value = "'6.01"
$ws.Range("E26").Value = "  -24.18%  "

$ws.Range("E27").Value = "  -23.05%  "

$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "'1.78"
$ws.Range("E28").Value = "  -18.38%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'6.68"
$ws.Range("E29").Value = "  -18.13%  "

$ws.Range("D30").Value = "'23.72"
$ws.Range("E30").Value = "  -19.25%  "

$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.37%  "

$ws.Range("E32").Value = "  -11.51%  "

$ws.Range("D33").Value = "'2.25"
$ws.Range("E33").Value = "  -16.10%  "

$ws.Range("D34").Value = "'50.12"
$ws.Range("E34").Value = "  -6.15%  "

$ws.Range("D35").Value = "'456.77"
$ws.Range("E35").Value = "  -16.71%  "

$ws.Range("D36").Value = "'5.23"
$ws.Range("E36").Value = "  -20.06%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.0376"
$ws.Range("E37").Value = "  -12.11%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.48"
$ws.Range("E38").Value = "  -21.52%  "

$ws.Range("D39").Value = "'0.112"
$ws.Range("E39").Value = "  -10.99%  "

$ws.Range("D40").Value = "'0.0719"
$ws.Range("E40").Value = "  -16.57%  "

$ws.Range("D41").Value = "'7.46"
$ws.Range("E41").Value = "  -20.37%  "

$ws.Range("D42").Value = "'2.552.45"
$ws.Range("E42").Value = "  -12.76%  "

$ws.Range("D43").Value = "'0.997"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("D44").Value = "'2.03"
$ws.Range("E44").Value = "  -24.33%  "

$ws.Range("E45").Value = "  -20.72%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'107.44"
$ws.Range("E46").Value = "  -11.11%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.0973"
$ws.Range("E47").Value = "  -14.96%  "

$ws.Range("D48").Value = "'1.70"
$ws.Range("E48").Value = "  -20.83%  "

$ws.Range("D49").Value = "'1.15"
$ws.Range("E49").Value = "  -7.12%  "

$ws.Range("D50").Value = "'0.0₃0439"
$ws.Range("E50").Value = "  -24.94%  "

$ws.Range("D51").Value = "'20.50"
$ws.Range("E51").Value = "  -23.15%  "
